$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.77954133333333
$ws.Range("H2").Value = 68.338624
$ws.Range("I2").Value = 0.8649343844704168
$ws.Range("J2").Value = 0.8649343844704167
$ws.Range("M2").Value = 70.46291600000001
$ws.Range("N2").Value = 211.388748
$ws.Range("O2").Value = 0.5276750397950939
$ws.Range("P2").Value = 0.5276750397950939
$ws.Range("Q2").Value = 1605.112907489195
$ws.Range("R2").Value = 14446.01616740275
$ws.Range("S2").Value = 0.4564042857455722
$ws.Range("T2").Value = 0.4564042857455722
$ws.Range("G3").Value = 22.77954133333333
$ws.Range("H3").Value = 68.338624
$ws.Range("I3").Value = 0.8649343844704168
$ws.Range("J3").Value = 0.8649343844704167
$ws.Range("O3").Value = 0.07361176802536967
$ws.Range("P3").Value = 0.07361176802536967
$ws.Range("Q3").Value = 223.9165965600782
$ws.Range("R3").Value = 2015.249369040704
$ws.Range("S3").Value = 0.06366934926680223
$ws.Range("T3").Value = 0.06366934926680222
$ws.Range("G4").Value = 22.77954133333333
$ws.Range("H4").Value = 68.338624
$ws.Range("I4").Value = 0.8649343844704168
$ws.Range("J4").Value = 0.8649343844704167
$ws.Range("M4").Value = 42.505498
$ws.Range("N4").Value = 127.516494
$ws.Range("O4").Value = 0.3183105613832428
$ws.Range("P4").Value = 0.3183105613832428
$ws.Range("Q4").Value = 968.2557485849171
$ws.Range("R4").Value = 8714.301737264255
$ws.Range("S4").Value = 0.275317749480448
$ws.Range("T4").Value = 0.2753177494804479
$ws.Range("G5").Value = 22.77954133333333
$ws.Range("H5").Value = 68.338624
$ws.Range("I5").Value = 0.8649343844704168
$ws.Range("J5").Value = 0.8649343844704167
$ws.Range("M5").Value = 10.73653933333333
$ws.Range("N5").Value = 32.209618
$ws.Range("O5").Value = 0.08040263079629371
$ws.Range("P5").Value = 0.08040263079629371
$ws.Range("Q5").Value = 244.5734415206257
$ws.Range("R5").Value = 2201.160973685632
$ws.Range("S5").Value = 0.06954299997759449
$ws.Range("T5").Value = 0.06954299997759447
$ws.Range("I6").Value = 0.008798055815159926
$ws.Range("J6").Value = 0.008798055815159925
$ws.Range("M6").Value = 70.46291600000001
$ws.Range("N6").Value = 211.388748
$ws.Range("O6").Value = 0.5276750397950939
$ws.Range("P6").Value = 0.5276750397950939
$ws.Range("Q6").Value = 16.327103192192
$ws.Range("R6").Value = 146.943928729728
$ws.Range("S6").Value = 0.004642514452383972
$ws.Range("T6").Value = 0.004642514452383971
$ws.Range("I7").Value = 0.008798055815159926
$ws.Range("J7").Value = 0.008798055815159925
$ws.Range("O7").Value = 0.07361176802536967
$ws.Range("P7").Value = 0.07361176802536967
$ws.Range("S7").Value = 0.0006476404437398072
$ws.Range("T7").Value = 0.0006476404437398071
$ws.Range("I8").Value = 0.008798055815159926
$ws.Range("J8").Value = 0.008798055815159925
$ws.Range("M8").Value = 42.505498
$ws.Range("N8").Value = 127.516494
$ws.Range("O8").Value = 0.3183105613832428
$ws.Range("P8").Value = 0.3183105613832428
$ws.Range("Q8").Value = 9.849033952575999
$ws.Range("R8").Value = 88.641305573184
$ws.Range("S8").Value = 0.00280051408560466
$ws.Range("T8").Value = 0.00280051408560466
$ws.Range("I9").Value = 0.008798055815159926
$ws.Range("J9").Value = 0.008798055815159925
$ws.Range("M9").Value = 10.73653933333333
$ws.Range("N9").Value = 32.209618
$ws.Range("O9").Value = 0.08040263079629371
$ws.Range("P9").Value = 0.08040263079629371
$ws.Range("Q9").Value = 2.487785002005333
$ws.Range("R9").Value = 22.390065018048
$ws.Range("S9").Value = 0.0007073868334314885
$ws.Range("T9").Value = 0.0007073868334314884
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.07785033333333334
$ws.Range("H10").Value = 0.233551
$ws.Range("I10").Value = 0.002955960752552617
$ws.Range("J10").Value = 0.002955960752552617
$ws.Range("M10").Value = 70.46291600000001
$ws.Range("N10").Value = 211.388748
$ws.Range("O10").Value = 0.5276750397950939
$ws.Range("P10").Value = 0.5276750397950939
$ws.Range("Q10").Value = 5.485561498238668
$ws.Range("R10").Value = 49.37005348414801
$ws.Range("S10").Value = 0.001559786707735938
$ws.Range("T10").Value = 0.001559786707735938
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.07785033333333334
$ws.Range("H11").Value = 0.233551
$ws.Range("I11").Value = 0.002955960752552617
$ws.Range("J11").Value = 0.002955960752552617
$ws.Range("O11").Value = 0.07361176802536967
$ws.Range("P11").Value = 0.07361176802536967
$ws.Range("Q11").Value = 0.7652472640245556
$ws.Range("R11").Value = 6.887225376221
$ws.Range("S11").Value = 0.0002175934972090004
$ws.Range("T11").Value = 0.0002175934972090004
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.07785033333333334
$ws.Range("H12").Value = 0.233551
$ws.Range("I12").Value = 0.002955960752552617
$ws.Range("J12").Value = 0.002955960752552617
$ws.Range("M12").Value = 42.505498
$ws.Range("N12").Value = 127.516494
$ws.Range("O12").Value = 0.3183105613832428
$ws.Range("P12").Value = 0.3183105613832428
$ws.Range("Q12").Value = 3.309067187799333
$ws.Range("R12").Value = 29.781604690194
$ws.Range("S12").Value = 0.0009409135265718565
$ws.Range("T12").Value = 0.0009409135265718564
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.07785033333333334
$ws.Range("H13").Value = 0.233551
$ws.Range("I13").Value = 0.002955960752552617
$ws.Range("J13").Value = 0.002955960752552617
$ws.Range("M13").Value = 10.73653933333333
$ws.Range("N13").Value = 32.209618
$ws.Range("O13").Value = 0.08040263079629371
$ws.Range("P13").Value = 0.08040263079629371
$ws.Range("Q13").Value = 0.8358431659464445
$ws.Range("R13").Value = 7.522588493518
$ws.Range("S13").Value = 0.0002376670210358226
$ws.Range("T13").Value = 0.0002376670210358226
$ws.Range("G14").Value = 3.247624
$ws.Range("H14").Value = 9.742872
$ws.Range("I14").Value = 0.1233115989618705
$ws.Range("J14").Value = 0.1233115989618705
$ws.Range("M14").Value = 70.46291600000001
$ws.Range("N14").Value = 211.388748
$ws.Range("O14").Value = 0.5276750397950939
$ws.Range("P14").Value = 0.5276750397950939
$ws.Range("Q14").Value = 228.837057111584
$ws.Range("R14").Value = 2059.533514004256
$ws.Range("S14").Value = 0.06506845288940169
$ws.Range("T14").Value = 0.06506845288940168
$ws.Range("G15").Value = 3.247624
$ws.Range("H15").Value = 9.742872
$ws.Range("I15").Value = 0.1233115989618705
$ws.Range("J15").Value = 0.1233115989618705
$ws.Range("O15").Value = 0.07361176802536967
$ws.Range("P15").Value = 0.07361176802536967
$ws.Range("Q15").Value = 31.92324649323466
$ws.Range("R15").Value = 287.309218439112
$ws.Range("S15").Value = 0.009077184817618629
$ws.Range("T15").Value = 0.009077184817618628
$ws.Range("G16").Value = 3.247624
$ws.Range("H16").Value = 9.742872
$ws.Range("I16").Value = 0.1233115989618705
$ws.Range("J16").Value = 0.1233115989618705
$ws.Range("M16").Value = 42.505498
$ws.Range("N16").Value = 127.516494
$ws.Range("O16").Value = 0.3183105613832428
$ws.Range("P16").Value = 0.3183105613832428
$ws.Range("Q16").Value = 138.041875436752
$ws.Range("R16").Value = 1242.376878930768
$ws.Range("S16").Value = 0.03925138429061831
$ws.Range("T16").Value = 0.0392513842906183
$ws.Range("G17").Value = 3.247624
$ws.Range("H17").Value = 9.742872
$ws.Range("I17").Value = 0.1233115989618705
$ws.Range("J17").Value = 0.1233115989618705
$ws.Range("M17").Value = 10.73653933333333
$ws.Range("N17").Value = 32.209618
$ws.Range("O17").Value = 0.08040263079629371
$ws.Range("P17").Value = 0.08040263079629371
$ws.Range("Q17").Value = 34.86824281587733
$ws.Range("R17").Value = 313.814185342896
$ws.Range("S17").Value = 0.009914576964231911
$ws.Range("T17").Value = 0.009914576964231909
